# Updates the cryptos price table: Price (D) and Volume(1h) (E) columns,
# plus one coin replacement (row 51: EnergySwap -> USDD) per the scraped refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.660.18"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.595.01"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'210.84"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").Value = "'19.43"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").Value = "'0.0841"
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").Value = "1.819.36"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "1.618.31"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Value = "'0.520"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "'64.42"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").Value = "26.649.87"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D20").Value = "'207.51"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").Value = "'6.78"
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").Value = "'4.22"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").Value = "'2.32"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").Value = "'145.49"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'7.18"
$ws.Range("E27").Value = "  -2.37%  "
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").Value = "'15.21"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").Value = "'0.0504"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("D31").Value = "'1.15"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").Value = "'0.662"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "'2.92"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").Value = "1.278.81"
$ws.Range("E35").Value = "  -3.87%  "
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").Value = "'0.838"
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("D43").Value = "'0.784"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "1.732.03"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").Value = "'0.904"
$ws.Range("E46").Value = "  +8.90%  "
$ws.Range("D47").Value = "'89.95"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.04%  "
